$d = $word.ActiveDocument

# --- 1 & 2. Replace text content in existing paragraphs (by paragraph index), using
#    direct Range.Text assignment to avoid AutoFormat/AutoCorrect mangling (e.g. smart quotes).

$p0 = $d.Paragraphs(3)
$check0 = $p0.Range.Text.TrimEnd([char]13)
if ($check0 -ne "JÃO, Jema") { throw "Paragraph 3 did not contain expected text for replacement 0: [$check0]" }
$p0.Range.Text = "BAY, jhony"

$p1 = $d.Paragraphs(6)
$check1 = $p1.Range.Text.TrimEnd([char]13)
if ($check1 -ne "O tabaco, uma planta do gênero Nicotiana, é cultivado mundialmente e amplamente consumido, principalmente sob a forma de cigarros, charutos e rapé. Este artigo tem como objetivo analisar os diversos aspectos relacionados ao tabaco, desde sua história e cultivo até seus efeitos na saúde pública. A pesquisa aborda a composição química do tabaco, com ênfase na nicotina, o principal componente psicoativo responsável pela dependência. A revisão de literatura explora os inúmeros estudos que comprovam a relação entre o consumo de tabaco e diversas doenças, incluindo câncer, doenças cardiovasculares e respiratórias. Além disso, são discutidas as políticas públicas de controle do tabagismo, como a proibição da publicidade, o aumento de impostos e as advertências nas embalagens, que visam reduzir o consumo e proteger a saúde da população. A metodologia utilizada consiste em uma revisão bibliográfica abrangente, analisando artigos científicos, relatórios de organizações de saúde e documentos governamentais. Os resultados esperados apontam para a necessidade de fortalecer as ações de prevenção e tratamento do tabagismo, bem como de promover a educação sobre os riscos associados ao consumo de tabaco. A conclusão ressalta a importância de um esforço conjunto entre governos, profissionais de saúde e sociedade civil para combater essa epidemia global e proteger as futuras gerações dos danos causados pelo tabaco.") { throw "Paragraph 6 did not contain expected text for replacement 1: [$check1]" }
$p1.Range.Text = "Este artigo investiga a relevância econômica da cultura do tabaco na região Sul do Brasil, analisando seu impacto no desenvolvimento socioeconômico dos estados do Paraná, Santa Catarina e Rio Grande do Sul. A produção de tabaco, embora controversa devido aos seus efeitos na saúde, historicamente representa uma parcela significativa da economia da região, gerando empregos, renda e divisas através da exportação. O estudo busca compreender a dinâmica da cadeia produtiva do tabaco, desde o cultivo até a industrialização e comercialização, identificando os principais atores envolvidos e os desafios enfrentados pelo setor. A pesquisa também examina as políticas públicas implementadas para regular a produção e o consumo de tabaco, bem como as iniciativas voltadas para a diversificação da economia nas regiões produtoras. Além disso, o artigo aborda os impactos sociais e ambientais da cultura do tabaco, considerando a saúde dos trabalhadores rurais, a utilização de agrotóxicos e o desmatamento. Através de uma revisão bibliográfica abrangente e da análise de dados estatísticos, o trabalho busca fornecer uma visão abrangente da importância do tabaco para a economia do Sul do Brasil, bem como os desafios e oportunidades para um futuro mais sustentável e diversificado."

$p2 = $d.Paragraphs(8)
$check2 = $p2.Range.Text.TrimEnd([char]13)
if ($check2 -ne "Tabaco; Nicotina; Tabagismo; Saúde Pública; Políticas de Controle.") { throw "Paragraph 8 did not contain expected text for replacement 2: [$check2]" }
$p2.Range.Text = "Tabaco; Economia; Região Sul; Desenvolvimento; Agricultura."

$p3 = $d.Paragraphs(10)
$check3 = $p3.Range.Text.TrimEnd([char]13)
if ($check3 -ne "Tobacco, a plant of the genus Nicotiana, is cultivated worldwide and widely consumed, mainly in the form of cigarettes, cigars, and snuff. This article aims to analyze the various aspects related to tobacco, from its history and cultivation to its effects on public health. The research addresses the chemical composition of tobacco, with emphasis on nicotine, the main psychoactive component responsible for addiction. The literature review explores numerous studies that prove the relationship between tobacco consumption and various diseases, including cancer, cardiovascular, and respiratory diseases. In addition, public policies to control smoking are discussed, such as the prohibition of advertising, increased taxes, and warnings on packaging, which aim to reduce consumption and protect public health. The methodology used consists of a comprehensive literature review, analyzing scientific articles, reports from health organizations, and government documents. The expected results point to the need to strengthen prevention and treatment actions for smoking, as well as to promote education about the risks associated with tobacco consumption. The conclusion highlights the importance of a joint effort between governments, health professionals, and civil society to combat this global epidemic and protect future generations from the harm caused by tobacco.") { throw "Paragraph 10 did not contain expected text for replacement 3: [$check3]" }
$p3.Range.Text = "This article investigates the economic relevance of tobacco cultivation in the Southern region of Brazil, analyzing its impact on the socioeconomic development of the states of Paraná, Santa Catarina, and Rio Grande do Sul. Tobacco production, although controversial due to its health effects, has historically represented a significant portion of the region's economy, generating jobs, income, and foreign exchange through exports. The study seeks to understand the dynamics of the tobacco production chain, from cultivation to industrialization and commercialization, identifying the main actors involved and the challenges faced by the sector. The research also examines the public policies implemented to regulate tobacco production and consumption, as well as initiatives aimed at diversifying the economy in producing regions. Furthermore, the article addresses the social and environmental impacts of tobacco cultivation, considering the health of rural workers, the use of pesticides, and deforestation. Through a comprehensive literature review and analysis of statistical data, the work seeks to provide a comprehensive view of the importance of tobacco for the economy of Southern Brazil, as well as the challenges and opportunities for a more sustainable and diversified future."

$p4 = $d.Paragraphs(12)
$check4 = $p4.Range.Text.TrimEnd([char]13)
if ($check4 -ne "Tobacco; Nicotine; Smoking; Public Health; Control Policies.") { throw "Paragraph 12 did not contain expected text for replacement 4: [$check4]" }
$p4.Range.Text = "Tobacco; Economy; Southern Region; Development; Agriculture."

$p5 = $d.Paragraphs(14)
$check5 = $p5.Range.Text.TrimEnd([char]13)
if ($check5 -ne "O tabaco, planta originária das Américas e pertencente ao gênero Nicotiana, possui uma longa história de uso, tanto em rituais religiosos quanto em práticas medicinais. Com a colonização, seu consumo se espalhou pelo mundo, transformando-se em um hábito socialmente aceito e, posteriormente, em um problema de saúde pública global. A popularização do cigarro, em particular, intensificou a dependência e aumentou exponencialmente a exposição da população aos seus componentes tóxicos.") { throw "Paragraph 14 did not contain expected text for replacement 5: [$check5]" }
$p5.Range.Text = "A cultura do tabaco possui uma longa história na região Sul do Brasil, remontando ao período colonial. Ao longo dos séculos, o tabaco se consolidou como uma importante atividade econômica, influenciando o desenvolvimento social e econômico dos estados do Paraná, Santa Catarina e Rio Grande do Sul. A produção de tabaco gera empregos diretos e indiretos, além de contribuir para a arrecadação de impostos e a geração de divisas através da exportação. No entanto, a produção e o consumo de tabaco são temas controversos devido aos seus efeitos prejudiciais à saúde."

$p6 = $d.Paragraphs(15)
$check6 = $p6.Range.Text.TrimEnd([char]13)
if ($check6 -ne "A presente pesquisa se justifica pela magnitude do problema do tabagismo e seus impactos devastadores na saúde individual e coletiva. A Organização Mundial da Saúde (OMS) estima que o tabaco seja responsável por mais de 8 milhões de mortes por ano em todo o mundo, tornando-se uma das principais causas de morbidade e mortalidade evitáveis. O alto custo do tratamento das doenças relacionadas ao tabaco sobrecarrega os sistemas de saúde, gerando um impacto econômico significativo.") { throw "Paragraph 15 did not contain expected text for replacement 6: [$check6]" }
$p6.Range.Text = "A relevância da pesquisa reside na necessidade de compreender a complexa relação entre o tabaco e a economia da região Sul, considerando os impactos positivos e negativos da cultura. Apesar dos esforços para diversificar a economia local, o tabaco ainda representa uma fonte significativa de renda para muitos produtores rurais, principalmente em pequenas propriedades."

$p7 = $d.Paragraphs(16)
$check7 = $p7.Range.Text.TrimEnd([char]13)
if ($check7 -ne "O problema central desta pesquisa reside na necessidade de compreender os mecanismos de dependência do tabaco, os fatores que contribuem para o início e a manutenção do tabagismo, bem como os efeitos nocivos do tabaco na saúde humana. Além disso, busca-se analisar a eficácia das políticas públicas de controle do tabagismo e identificar estratégias inovadoras para reduzir o consumo e proteger a população.") { throw "Paragraph 16 did not contain expected text for replacement 7: [$check7]" }
$p7.Range.Text = "O problema de pesquisa que se busca responder é: qual a real importância econômica da cultura do tabaco para a região Sul do Brasil, considerando os seus impactos sociais, ambientais e de saúde pública?"

$p8 = $d.Paragraphs(17)
$check8 = $p8.Range.Text.TrimEnd([char]13)
if ($check8 -ne "O objetivo geral deste artigo é analisar os diversos aspectos relacionados ao tabaco, desde sua composição química e efeitos na saúde até as políticas de controle do tabagismo. Especificamente, pretende-se: (1) revisar a literatura científica sobre os efeitos do tabaco na saúde; (2) analisar as políticas públicas de controle do tabagismo implementadas em diferentes países; (3) discutir as estratégias de prevenção e tratamento do tabagismo; e (4) identificar as lacunas no conhecimento e as necessidades de pesquisa futura.") { throw "Paragraph 17 did not contain expected text for replacement 8: [$check8]" }
$p8.Range.Text = "O objetivo geral desta pesquisa é analisar a importância econômica do tabaco na região Sul do Brasil, identificando os seus principais impactos e desafios. Os objetivos específicos incluem: a) descrever a cadeia produtiva do tabaco na região Sul; b) analisar os impactos sociais e ambientais da cultura do tabaco; c) examinar as políticas públicas voltadas para a regulação do tabaco e a diversificação da economia; d) identificar os desafios e oportunidades para um futuro mais sustentável na produção de tabaco."

$p9 = $d.Paragraphs(19)
$check9 = $p9.Range.Text.TrimEnd([char]13)
if ($check9 -ne "O tabaco contém milhares de substâncias químicas, sendo a nicotina o principal componente psicoativo responsável pela dependência. A nicotina atua no sistema nervoso central, liberando neurotransmissores como a dopamina, que promovem sensações de prazer e recompensa. A exposição repetida à nicotina leva à tolerância e à necessidade de doses cada vez maiores para obter os mesmos efeitos, caracterizando a dependência.") { throw "Paragraph 19 did not contain expected text for replacement 9: [$check9]" }
$p9.Range.Text = "A importância econômica do tabaco na região Sul do Brasil é amplamente reconhecida na literatura. Segundo dados da Associação dos Fumicultores do Brasil (Afubra), a região Sul concentra a maior parte da produção nacional de tabaco, sendo responsável por cerca de 95% da produção total. A cultura do tabaco gera empregos e renda para milhares de famílias rurais, além de contribuir para a arrecadação de impostos e a geração de divisas através da exportação."

$p10 = $d.Paragraphs(20)
$check10 = $p10.Range.Text.TrimEnd([char]13)
if ($check10 -ne "Os efeitos nocivos do tabaco na saúde são amplamente documentados. O consumo de tabaco está associado a um risco aumentado de diversas doenças, incluindo câncer de pulmão, boca, garganta, esôfago, bexiga, pâncreas, rim e colo do útero. Além disso, o tabaco contribui para o desenvolvimento de doenças cardiovasculares, como infarto do miocárdio e acidente vascular cerebral, e doenças respiratórias, como bronquite crônica e enfisema pulmonar.") { throw "Paragraph 20 did not contain expected text for replacement 10: [$check10]" }
$p10.Range.Text = "No entanto, a produção de tabaco também apresenta desafios significativos, incluindo os impactos na saúde dos trabalhadores rurais, a utilização de agrotóxicos e o desmatamento. Segundo estudo de Rosemberg (2003, p. 25), `"o uso intensivo de agrotóxicos na cultura do tabaco representa um grave problema de saúde pública, expondo os trabalhadores rurais a riscos de intoxicação e doenças crônicas.`""

$p11 = $d.Paragraphs(21)
$check11 = $p11.Range.Text.TrimEnd([char]13)
if ($check11 -ne "As políticas públicas de controle do tabagismo têm se mostrado eficazes na redução do consumo e na proteção da saúde da população. A proibição da publicidade de cigarros, o aumento de impostos sobre o tabaco e as advertências nas embalagens são algumas das medidas que têm contribuído para diminuir a prevalência do tabagismo em diversos países. Segundo Eriksen et al. (2015, p. 36), `"as políticas abrangentes de controle do tabagismo, que incluem medidas como o aumento de impostos, a proibição da publicidade e as advertências nas embalagens, são eficazes na redução do consumo e na proteção da saúde da população`". Além disso, a criação de ambientes livres de fumo e a oferta de tratamento para a dependência de nicotina também são importantes para auxiliar os fumantes a abandonarem o hábito. De acordo com Prochaska et al. (2013, p. 1428), `"o tratamento da dependência de nicotina, que inclui aconselhamento e medicamentos, é eficaz e custo-efetivo`".") { throw "Paragraph 21 did not contain expected text for replacement 11: [$check11]" }
$p11.Range.Text = "Além dos impactos na saúde e no meio ambiente, a produção de tabaco também enfrenta desafios relacionados à regulamentação do setor. A Convenção-Quadro para o Controle do Tabaco (CQCT), um tratado internacional da Organização Mundial da Saúde (OMS), estabelece medidas para reduzir o consumo de tabaco, incluindo o aumento de impostos, a proibição da publicidade e a adoção de embalagens padronizadas. Segundo Iglesias (2006, p. 112), `"a implementação da CQCT representa um desafio para a indústria do tabaco, que busca resistir às medidas de controle do tabagismo.`" A diversificação da economia nas regiões produtoras de tabaco é fundamental para reduzir a dependência da cultura e promover um desenvolvimento mais sustentável."

$p12 = $d.Paragraphs(23)
$check12 = $p12.Range.Text.TrimEnd([char]13)
if ($check12 -ne "Este artigo foi desenvolvido a partir de uma revisão bibliográfica abrangente, com abordagem qualitativa, sobre o tema do tabaco. A pesquisa foi conduzida em bases de dados científicas como PubMed, Scopus e Web of Science, utilizando os seguintes termos de busca: `"tabaco`", `"nicotina`", `"tabagismo`", `"saúde pública`" e `"políticas de controle`". Além disso, foram consultados relatórios de organizações de saúde, como a Organização Mundial da Saúde (OMS) e o Instituto Nacional de Câncer (INCA), bem como documentos governamentais relacionados ao controle do tabagismo.") { throw "Paragraph 23 did not contain expected text for replacement 12: [$check12]" }
$p12.Range.Text = "A presente pesquisa adotou uma abordagem metodológica mista, combinando métodos qualitativos e quantitativos. A pesquisa bibliográfica foi realizada em bases de dados acadêmicas, como Scielo, Google Scholar e periódicos especializados, com o objetivo de identificar os principais estudos e artigos sobre a importância econômica do tabaco na região Sul do Brasil."

$p13 = $d.Paragraphs(24)
$check13 = $p13.Range.Text.TrimEnd([char]13)
if ($check13 -ne "Os artigos e documentos selecionados foram analisados criticamente, buscando identificar os principais conceitos, evidências empíricas e lacunas no conhecimento. A análise dos dados foi realizada de forma descritiva e interpretativa, com o objetivo de sintetizar as informações relevantes e apresentar uma visão abrangente sobre o tema do tabaco.") { throw "Paragraph 24 did not contain expected text for replacement 13: [$check13]" }
$p13.Range.Text = "Além da pesquisa bibliográfica, foram utilizados dados estatísticos da Associação dos Fumicultores do Brasil (Afubra), do Instituto Brasileiro de Geografia e Estatística (IBGE) e do Ministério da Agricultura, Pecuária e Abastecimento (MAPA) para analisar a produção, a comercialização e a exportação de tabaco na região Sul. Os dados foram analisados de forma descritiva, utilizando tabelas e gráficos para apresentar os principais resultados."

$p14 = $d.Paragraphs(25)
$check14 = $p14.Range.Text.TrimEnd([char]13)
if ($check14 -ne "A pesquisa bibliográfica permitiu identificar os principais fatores de risco associados ao consumo de tabaco, as estratégias de prevenção e tratamento do tabagismo mais eficazes e as políticas públicas de controle do tabagismo que têm se mostrado bem-sucedidas em diferentes contextos. A análise dos dados também possibilitou identificar as áreas em que são necessárias mais pesquisas, como a avaliação dos efeitos de novas formas de consumo de tabaco (como os cigarros eletrônicos) e o desenvolvimento de intervenções mais eficazes para grupos específicos de fumantes.") { throw "Paragraph 25 did not contain expected text for replacement 14: [$check14]" }
$p14.Range.Text = "A abordagem qualitativa foi utilizada para analisar as políticas públicas voltadas para a regulação do tabaco e a diversificação da economia. Foram analisados documentos oficiais, como leis, decretos e portarias, com o objetivo de identificar as principais medidas implementadas para controlar o tabagismo e promover o desenvolvimento de atividades alternativas."

$p15 = $d.Paragraphs(27)
$check15 = $p15.Range.Text.TrimEnd([char]13)
if ($check15 -ne "Espera-se que esta revisão bibliográfica demonstre a magnitude do problema do tabagismo e seus impactos na saúde pública. Os resultados esperados apontam para a necessidade de fortalecer as ações de prevenção e tratamento do tabagismo, bem como de promover a educação sobre os riscos associados ao consumo de tabaco. A análise das políticas públicas de controle do tabagismo deverá revelar as estratégias mais eficazes para reduzir o consumo e proteger a população, como o aumento de impostos, a proibição da publicidade e as advertências nas embalagens.") { throw "Paragraph 27 did not contain expected text for replacement 15: [$check15]" }
$p15.Range.Text = "Os resultados da pesquisa indicam que a cultura do tabaco ainda representa uma importante fonte de renda para muitos produtores rurais na região Sul do Brasil, especialmente em pequenas propriedades. No entanto, a importância econômica do tabaco tem diminuído nos últimos anos, devido à crescente conscientização sobre os seus efeitos prejudiciais à saúde e às políticas públicas de controle do tabagismo."

$p16 = $d.Paragraphs(28)
$check16 = $p16.Range.Text.TrimEnd([char]13)
if ($check16 -ne "A discussão dos resultados deverá abordar os desafios enfrentados no controle do tabagismo, como a resistência da indústria do tabaco, a falta de recursos para implementar políticas eficazes e a dificuldade em alcançar grupos específicos de fumantes, como os jovens e as populações de baixa renda. Além disso, a discussão deverá explorar as oportunidades para inovar nas estratégias de prevenção e tratamento do tabagismo, como o uso de tecnologias digitais e a personalização das intervenções.") { throw "Paragraph 28 did not contain expected text for replacement 16: [$check16]" }
$p16.Range.Text = "A cadeia produtiva do tabaco na região Sul é caracterizada pela presença de grandes empresas fumageiras, que controlam a produção, a comercialização e a exportação do produto. Os produtores rurais, em geral, são dependentes das empresas fumageiras, que fornecem insumos e assistência técnica e compram a produção a preços pré-definidos."

$p17 = $d.Paragraphs(29)
$check17 = $p17.Range.Text.TrimEnd([char]13)
if ($check17 -ne "Espera-se que a análise dos dados revele a importância de um esforço conjunto entre governos, profissionais de saúde e sociedade civil para combater essa epidemia global e proteger as futuras gerações dos danos causados pelo tabaco. A discussão dos resultados deverá destacar a necessidade de investir em pesquisa para aprimorar as estratégias de prevenção e tratamento do tabagismo, bem como para monitorar os efeitos das políticas públicas de controle do tabagismo.") { throw "Paragraph 29 did not contain expected text for replacement 17: [$check17]" }
$p17.Range.Text = "Os impactos sociais e ambientais da cultura do tabaco são significativos. A utilização intensiva de agrotóxicos na cultura do tabaco representa um grave problema de saúde pública, expondo os trabalhadores rurais a riscos de intoxicação e doenças crônicas. Além disso, a cultura do tabaco contribui para o desmatamento, especialmente em áreas de Mata Atlântica."

$p18 = $d.Paragraphs(31)
$check18 = $p18.Range.Text.TrimEnd([char]13)
if ($check18 -ne "Este artigo teve como objetivo analisar os diversos aspectos relacionados ao tabaco, desde sua composição química e efeitos na saúde até as políticas de controle do tabagismo. Através da revisão da literatura científica, foi possível constatar a magnitude do problema do tabagismo e seus impactos devastadores na saúde individual e coletiva. As políticas públicas de controle do tabagismo têm se mostrado eficazes na redução do consumo, mas ainda há muitos desafios a serem superados.") { throw "Paragraph 31 did not contain expected text for replacement 18: [$check18]" }
$p18.Range.Text = "Este estudo investigou a importância econômica do tabaco na região Sul do Brasil, analisando os seus impactos sociais, ambientais e de saúde pública. Os resultados indicam que a cultura do tabaco ainda representa uma importante fonte de renda para muitos produtores rurais, mas a sua importância econômica tem diminuído nos últimos anos, devido à crescente conscientização sobre os seus efeitos prejudiciais à saúde e às políticas públicas de controle do tabagismo."

$p19 = $d.Paragraphs(32)
$check19 = $p19.Range.Text.TrimEnd([char]13)
if ($check19 -ne "As contribuições deste artigo residem na síntese das informações relevantes sobre o tabaco e na identificação das lacunas no conhecimento e das necessidades de pesquisa futura. A pesquisa sugere que é fundamental fortalecer as ações de prevenção e tratamento do tabagismo, promover a educação sobre os riscos associados ao consumo de tabaco e investir em pesquisa para aprimorar as estratégias de controle do tabagismo.") { throw "Paragraph 32 did not contain expected text for replacement 19: [$check19]" }
$p19.Range.Text = "A diversificação da economia nas regiões produtoras de tabaco é fundamental para reduzir a dependência da cultura e promover um desenvolvimento mais sustentável. As políticas públicas devem incentivar a diversificação, oferecendo apoio técnico e financeiro aos produtores rurais para que possam desenvolver atividades alternativas."

$p20 = $d.Paragraphs(33)
$check20 = $p20.Range.Text.TrimEnd([char]13)
if ($check20 -ne "Como trabalhos futuros, sugere-se a realização de estudos que avaliem os efeitos de novas formas de consumo de tabaco (como os cigarros eletrônicos) e desenvolvam intervenções mais eficazes para grupos específicos de fumantes. Além disso, é importante monitorar os efeitos das políticas públicas de controle do tabagismo e avaliar a eficácia de diferentes estratégias de prevenção e tratamento do tabagismo em diferentes contextos.") { throw "Paragraph 33 did not contain expected text for replacement 20: [$check20]" }
$p20.Range.Text = "Sugere-se para trabalhos futuros a realização de estudos de caso em regiões específicas da região Sul, com o objetivo de analisar em profundidade os impactos sociais e ambientais da cultura do tabaco e as alternativas de diversificação econômica."

$p21 = $d.Paragraphs(35)
$check21 = $p21.Range.Text.TrimEnd([char]13)
if ($check21 -ne "ERIKSEN, M.; MACKAY, J.; ROSS, H. The Tobacco Atlas. 5. ed. Atlanta, GA: American Cancer Society, 2015.") { throw "Paragraph 35 did not contain expected text for replacement 21: [$check21]" }
$p21.Range.Text = "IGLESIAS, R. O Mito da Liberdade: Tabaco, Propaganda e Direito. São Paulo: Martins Fontes, 2006."

$p22 = $d.Paragraphs(36)
$check22 = $p22.Range.Text.TrimEnd([char]13)
if ($check22 -ne "INCA - INSTITUTO NACIONAL DE CÂNCER JOSÉ ALENCAR GOMES DA SILVA. Tabagismo. Rio de Janeiro: INCA, 2021. Disponível em: [https://www.inca.gov.br/causas-e-prevencao/tabagismo](https://www.inca.gov.br/causas-e-prevencao/tabagismo). Acesso em: 15 nov. 2023.") { throw "Paragraph 36 did not contain expected text for replacement 22: [$check22]" }
$p22.Range.Text = "ROSEMBERG, J. Tabagismo: sério problema de saúde pública no Brasil. São Paulo: Almed, 2003."

$p23 = $d.Paragraphs(37)
$check23 = $p23.Range.Text.TrimEnd([char]13)
if ($check23 -ne "PROCHASKA, J. J.; BENOWITZ, N. L. The past, present, and future of nicotine addiction therapy. Annual Review of Medicine, v. 64, p. 427-445, 2013.") { throw "Paragraph 37 did not contain expected text for replacement 23: [$check23]" }
$p23.Range.Text = "AFUBRA. Anuário Estatístico do Tabaco. Santa Cruz do Sul, RS: Afubra, 2023."

# --- 3. Insert two new paragraphs in the Resultados e Discussao section ---
# Paragraph 29 is the one ending "...Mata Atlantica." (formerly "Espera-se que a analise...")
$pInsertAfter = $d.Paragraphs(29)
$pInsertAfter.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs(30)
$newPara1.Range.Text = "As políticas públicas voltadas para a regulação do tabaco têm se intensificado nos últimos anos, com a implementação de medidas como o aumento de impostos, a proibição da publicidade e a adoção de embalagens padronizadas. No entanto, a implementação dessas políticas enfrenta resistência da indústria do tabaco, que busca defender os seus interesses econômicos."
$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs(31)
$newPara2.Range.Text = "A diversificação da economia nas regiões produtoras de tabaco é fundamental para reduzir a dependência da cultura e promover um desenvolvimento mais sustentável. Diversas iniciativas têm sido implementadas nesse sentido, incluindo o apoio à agricultura familiar, o fomento ao turismo rural e o desenvolvimento de atividades artesanais."

# --- 4. Apply justified alignment (jc=both) to body paragraphs ---
$jcParagraphIndices = @(6, 8, 10, 12, 14, 15, 16, 17, 19, 20, 21, 23, 24, 25, 27, 28, 29, 30, 31, 33, 34, 35)
foreach ($idx in $jcParagraphIndices) {
    $p = $d.Paragraphs($idx)
    $p.Range.ParagraphFormat.Alignment = 3  # wdAlignParagraphJustify
}

Write-Host "Done. Justified paragraph count: " $jcParagraphIndices.Count
